$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 620 ("「とても愛しています」...") entirely; rows below shift up by one.
$ws.Rows.Item(620).Delete()
